$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Updated GDP per Capita (Social Spending indicator source data) figures,
# and newly published Work Week (Togo) observations for 2009-2016.
$gdpByYear = [ordered]@{
    1950 = "870"
    1951 = "885"
    1952 = "899"
    1953 = "913"
    1954 = "928"
    1955 = "940"
    1956 = "955"
    1957 = "969"
    1958 = "982"
    1959 = "995"
    1960 = "1058"
    1961 = "1103"
    1962 = "1117"
    1963 = "1143"
    1964 = "1275"
    1965 = "1412"
    1966 = "1503"
    1967 = "1546"
    1968 = "1581"
    1969 = "1707"
    1970 = "1631"
    1971 = "1699"
    1972 = "1709"
    1973 = "1591"
    1974 = "1612"
    1975 = "1556"
    1976 = "1505"
    1977 = "1543"
    1978 = "1650"
    1979 = "1698"
    1980 = "1572"
    1981 = "1423"
    1982 = "1323"
    1983 = "1211"
    1984 = "1239"
    1985 = "1240"
    1986 = "1237"
    1987 = "1165"
    1988 = "1239"
    1989 = "1246"
    1990 = "1275"
    1991 = "1249.16794321487"
    1992 = "1181.87126506398"
    1993 = "1003.16487487302"
    1994 = "1145.77758007388"
    1995 = "1170.59771705596"
    1996 = "1215.73997885394"
    1997 = "1226.32704989982"
    1998 = "1174.85273698303"
    1999 = "1183.01296609916"
    2000 = "1150.36591070506"
    2001 = "1115.5502532021"
    2002 = "1081.74451524629"
    2003 = "1114.46381917196"
    2004 = "1121.80584889235"
    2005 = "1119.80426909117"
    2006 = "1144.8814458855"
    2007 = "1145.92064503615"
    2008 = "1150.35785293882"
    2009 = "1165.88883371318"
    2010 = "1191.14492660954"
    2011 = "1227"
    2012 = "1263"
    2013 = "1292"
    2014 = "1334"
    2015 = "1370"
    2016 = "1400"
}

$countryCode = 768
$countryName = "Togo"
$indicator = "GDP per Capita"
$firstDataRow = 2
$lastExistingRow = 60

$row = $firstDataRow
foreach ($year in $gdpByYear.Keys) {
    if ($row -gt $lastExistingRow) {
        # Year not present in the sheet yet - add a brand new row
        $ws.Cells.Item($row, 1).Value = $countryCode
        $ws.Cells.Item($row, 2).Value = $countryName
        $ws.Cells.Item($row, 3).Value = $indicator
        $ws.Cells.Item($row, 4).Value = $year
    }

    # Write the value as text (matching the source workbook's convention)
    $ws.Cells.Item($row, 5).Value = "'" + $gdpByYear[$year]
    $ws.Cells.Item($row, 5).ClearFormats()

    $row = $row + 1
}
